$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B3:J50").Value = 0
